$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.035.11'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '2.899.99'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.07'
$ws.Range('E5').Value = '  -3.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.55'
$ws.Range('E6').Value = '  -2.33%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '2.896.98'
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.96'
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.148'
$ws.Range('E11').Value = '  -1.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.431'
$ws.Range('E12').Value = '  -1.22%  '
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.43'
$ws.Range('E14').Value = '  -1.98%  '
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '3.378.42'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('D17').Value = '61.964.80'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.54'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '2.894.60'
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '427.67'
$ws.Range('E20').Value = '  -1.67%  '
$ws.Range('E21').Value = '  -3.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.656'
$ws.Range('E22').Value = '  -0.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.83'
$ws.Range('E23').Value = '  -1.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.40'
$ws.Range('E24').Value = '  -3.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.95'
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.87'
$ws.Range('E27').Value = '  -3.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.99'
$ws.Range('E28').Value = '  -4.37%  '
$ws.Range('E29').Value = '  +3.14%  '
$ws.Range('E30').Value = '  -4.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.47'
$ws.Range('E31').Value = '  -3.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.99'
$ws.Range('E32').Value = '  -5.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.60'
$ws.Range('E34').Value = '  -1.21%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.105'
$ws.Range('E35').Value = '  -3.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.950'
$ws.Range('E36').Value = '  -2.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.36'
$ws.Range('E37').Value = '  -2.80%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.78'
$ws.Range('E38').Value = '  -0.73%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.91'
$ws.Range('E39').Value = '  -5.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.88'
$ws.Range('E40').Value = '  -6.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.113'
$ws.Range('E41').Value = '  -2.22%  '
$ws.Range('B42').Value = 'Arweave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '40.76'
$ws.Range('E42').Value = '  +4.50%  '
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.07'
$ws.Range('E43').Value = '  -3.46%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.264'
$ws.Range('E44').Value = '  -3.20%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.703.87'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '133.32'
$ws.Range('E46').Value = '  -0.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0334'
$ws.Range('E47').Value = '  -1.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '346.42'
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000218'
$ws.Range('E50').Value = '  +13.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.102'
$ws.Range('E51').Value = '  -1.59%  '
